# Add a new "Cu SE (1ppb)" worksheet, placed right before "Cu Density",
# as a copy of the existing "SE (1ppb)" sheet but re-targeted to the Cu
# number density (F2/F11/F26 = 9.9899999999999993E+22), and make it the
# active/selected sheet (matching the author's edit).

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("SE (1ppb)")
$before = $wb.Worksheets.Item("Cu Density")

# Copy "SE (1ppb)" and drop the copy immediately before "Cu Density".
$source.Copy($before)

# Excel names the freshly-made copy "SE (1ppb) (2)"; rename + retarget it.
$newSheet = $wb.Worksheets.Item("SE (1ppb) (2)")
$newSheet.Name = "Cu SE (1ppb)"

# Cu overall number density (atoms/cm^3) -- same literal used in the three
# per-isotope blocks of the sheet (Th-232, U-238, U-235 rows).
$cuNumberDensity = 99899999999999992659968.0

$newSheet.Range("F2").Value = $cuNumberDensity
$newSheet.Range("F11").Value = $cuNumberDensity
$newSheet.Range("F26").Value = $cuNumberDensity

# Match the author's on-disk selection/active-tab state.
$newSheet.Select() | Out-Null
$newSheet.Range("F27").Select() | Out-Null
